# fix FY input range
# Clears the hard-coded FY period-end-date text labels that were typed
# directly into the "Data Entry" input column (A7:A16). These were plain
# text strings (e.g. "2015-12-31" ... "2024-09-30") instead of being left
# blank for the user to fill in, which fed downstream YEAR()/date math with
# bogus text. Clearing them lets the workbook's array formulas recompute
# cleanly from the (now empty) input range.

$wb = $excel.ActiveWorkbook

$dataEntry = $wb.Worksheets.Item("Data Entry")
$dataEntry.Range("A7:A16").ClearContents()

# Restore the UI selection state captured in the edit.
$dataEntry.Activate()
$dataEntry.Range("A21").Select()

$fcfData = $wb.Worksheets.Item("FCF DATA")
$fcfData.Application.ActiveWindow.ScrollRow = 40
$fcfData.Activate()
$fcfData.Range("G52").Select()

# Chart sheet zoom levels were bumped as part of this edit.
$freeCashFlowGraph = $wb.Sheets.Item("Free Cash Flow Graph")
$freeCashFlowGraph.Activate()
$excel.ActiveWindow.Zoom = 118

$growthYoYGraph = $wb.Sheets.Item("Growth YoY Graph")
$growthYoYGraph.Activate()
$excel.ActiveWindow.Zoom = 118

$dataEntry.Activate()
$excel.CalculateFullRebuild()
